$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 117.745847958593
$ws.Range("D2").Value = 401567.231247708
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 402818.6746646345
